$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: nvidia/NV-Embed-v2
$ws.Range("D2").Value = 92.31
$ws.Range("E2").Value = 98.45999999999999
$ws.Range("F2").Value = 98.45999999999999
$ws.Range("G2").Value = 98.45999999999999
$ws.Range("H2").Value = 98.45999999999999

# Row 3: dunzhang/stella_en_1.5B_v5
$ws.Range("D3").Value = 92.31

# Row 4: sentence-transformers/all-MiniLM-L6-v2
$ws.Range("D4").Value = 64.62
$ws.Range("E4").Value = 75.38
$ws.Range("F4").Value = 80
$ws.Range("G4").Value = 81.54000000000001
$ws.Range("H4").Value = 81.54000000000001

# Row 5: was MSMARCO -> now mixedbread-ai/mxbai-embed-large-v1
$ws.Range("A5").Value = "mixedbread-ai/mxbai-embed-large-v1"
$ws.Range("D5").Value = 87.69
$ws.Range("E5").Value = 98.45999999999999
$ws.Range("F5").Value = 98.45999999999999
$ws.Range("G5").Value = 98.45999999999999
$ws.Range("H5").Value = 98.45999999999999

# Row 6: was mxbai embed large -> now MSMARCO
$ws.Range("A6").Value = "MSMARCO"
$ws.Range("C6").Value = "{'instruction': None, 'query_instruction': None, 'load_in_8bit': False, 'trust_remote_code': False}"
$ws.Range("D6").Value = 76.92
$ws.Range("E6").Value = 86.15000000000001
$ws.Range("F6").Value = 87.69
$ws.Range("G6").Value = 89.23
$ws.Range("H6").Value = 89.23
